$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Transmitance values (column B, rows 3-18) from 1 to 100
# This adds the extra ordinary ray for polarimetric acquisitions
$ws.Range("B3:B18").Value = 100

# Move the active selection to B18
$ws.Range("B18").Select()
